$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-21 06:32:55"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
